$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "TC5"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1

$ws.Range("D6").Select()
